$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Final")
$ws.Name = "Data-Final"
